# Insert a new weekly record for "Poroto granado" (Macroferia Regional de Talca)
# at row 15, pushing the existing rows 15..108 down to 16..109.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 15..108 down by one row, creating a blank row 15.
$ws.Rows(15).Insert()

# Populate the new row 15 with the new weekly data point.
$ws.Range("A15").Value = 5
$ws.Range("B15").Value = "Macroferia Regional de Talca"
$ws.Range("C15").Value = "Maule"
$ws.Range("D15").Value = "2022-02-10"
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = 100112030
$ws.Range("G15").Value = "Poroto granado"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 20000
$ws.Range("L15").Value = 20000
$ws.Range("M15").Value = 20000
$ws.Range("N15").Value = "`$/saco 25 kilos"
$ws.Range("O15").Value = "Región del Maule"
$ws.Range("P15").Value = 800
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"
